$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1671
$ws.Range("I106").Value = 1099.25
$ws.Range("K106").Value = 1099.25
$ws.Range("M106").Value = -468.25

$ws.Range("H113").Value = 5892.5
$ws.Range("I113").Value = 4320
$ws.Range("J113").Value = 8251.25
$ws.Range("K113").Value = 4320
$ws.Range("L113").Value = 8251.25
$ws.Range("M113").Value = -1066
$ws.Range("N113").Value = -14759.25

$ws.Range("H132").Value = 4084494.8
$ws.Range("I132").Value = 4763959
$ws.Range("K132").Value = 14291877
$ws.Range("M132").Value = -14289347

$ws.Range("H135").Value = 487.08572
$ws.Range("I135").Value = 372.66666
$ws.Range("K135").Value = 3353.99994
$ws.Range("M135").Value = -818.9999399999997

$ws.Range("H137").Value = 1964059.2
$ws.Range("I137").Value = 2860614
$ws.Range("J137").Value = 2845.6875
$ws.Range("K137").Value = 8581842
$ws.Range("L137").Value = 8537.0625
$ws.Range("M137").Value = -8579292
$ws.Range("N137").Value = -13637.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2340.487
$ws.Range("I61").Value = 1797.2
$ws.Range("J61").Value = 2912.3684
$ws.Range("K61").Value = 1797.2
$ws.Range("L61").Value = 2912.3684
$ws.Range("M61").Value = -1585.2
$ws.Range("N61").Value = -3336.3684

$ws.Range("H86").Value = 23000
$ws.Range("I86").Value = 16000
$ws.Range("K86").Value = 16000
$ws.Range("M86").Value = -14814

$ws.Range("H89").Value = 23000
$ws.Range("I89").Value = 16000
$ws.Range("K89").Value = 48000
$ws.Range("M89").Value = -42072

$ws.Range("H110").Value = 1633
$ws.Range("I110").Value = 511.10526
$ws.Range("K110").Value = 511.10526
$ws.Range("M110").Value = 1533.89474

$ws.Range("H122").Value = 2068.3394
$ws.Range("I122").Value = 1702.6904
$ws.Range("J122").Value = 3165.2856
$ws.Range("K122").Value = 5108.0712
$ws.Range("L122").Value = 9495.856800000001
$ws.Range("M122").Value = -2658.0712
$ws.Range("N122").Value = -14395.8568

$ws.Range("H132").Value = 1768.9265
$ws.Range("I132").Value = 1595.5
$ws.Range("J132").Value = 2578.25
$ws.Range("K132").Value = 4786.5
$ws.Range("L132").Value = 7734.75
$ws.Range("M132").Value = -2256.5
$ws.Range("N132").Value = -12794.75

$ws.Range("H136").Value = 2340.487
$ws.Range("I136").Value = 1797.2
$ws.Range("J136").Value = 2912.3684
$ws.Range("K136").Value = 5391.6
$ws.Range("L136").Value = 8737.1052
$ws.Range("M136").Value = -2841.6
$ws.Range("N136").Value = -13837.1052

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3430.9
$ws.Range("I99").Value = 2976.125
$ws.Range("K99").Value = 2976.125
$ws.Range("M99").Value = -1478.125

$ws.Range("H107").Value = 1832.7273
$ws.Range("I107").Value = 954.2857
$ws.Range("J107").Value = 3370
$ws.Range("K107").Value = 954.2857
$ws.Range("L107").Value = 3370
$ws.Range("M107").Value = 965.7143
$ws.Range("N107").Value = -7210

$ws.Range("H135").Value = 28227.059
$ws.Range("J135").Value = 28178.75
$ws.Range("L135").Value = 28178.75
$ws.Range("N135").Value = -38318.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1033274.44
$ws.Range("I31").Value = 3032645.8
$ws.Range("J31").Value = 2348.625
$ws.Range("K31").Value = 3032645.8
$ws.Range("L31").Value = 2348.625
$ws.Range("M31").Value = -3032350.8
$ws.Range("N31").Value = -2938.625

$ws.Range("H34").Value = 1033274.44
$ws.Range("I34").Value = 3032645.8
$ws.Range("J34").Value = 2348.625
$ws.Range("K34").Value = 3032645.8
$ws.Range("L34").Value = 2348.625
$ws.Range("M34").Value = -3032443.8
$ws.Range("N34").Value = -2752.625

$ws.Range("H62").Value = 3569.7778
$ws.Range("I62").Value = 2359.5
$ws.Range("J62").Value = 3915.5715
$ws.Range("K62").Value = 2359.5
$ws.Range("L62").Value = 3915.5715
$ws.Range("M62").Value = -1735.5
$ws.Range("N62").Value = -5163.5715

$ws.Range("H65").Value = 3569.7778
$ws.Range("I65").Value = 2359.5
$ws.Range("J65").Value = 3915.5715
$ws.Range("K65").Value = 11797.5
$ws.Range("L65").Value = 19577.8575
$ws.Range("M65").Value = -8677.5
$ws.Range("N65").Value = -25817.8575

$ws.Range("H107").Value = 1072.1666
$ws.Range("I107").Value = 1095.6
$ws.Range("J107").Value = 1048.7333
$ws.Range("K107").Value = 1095.6
$ws.Range("L107").Value = 1048.7333
$ws.Range("M107").Value = 824.4000000000001
$ws.Range("N107").Value = -4888.7333

$ws.Range("H122").Value = 2261.652
$ws.Range("I122").Value = 1911.2632
$ws.Range("J122").Value = 3926
$ws.Range("K122").Value = 5733.7896
$ws.Range("L122").Value = 11778
$ws.Range("M122").Value = -3283.7896
$ws.Range("N122").Value = -16678

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2857926.8
$ws.Range("I113").Value = 12500480
$ws.Range("J113").Value = 873.85187
$ws.Range("K113").Value = 37501440
$ws.Range("L113").Value = 2621.55561
$ws.Range("M113").Value = -37499270
$ws.Range("N113").Value = -6961.555609999999

$ws.Range("H137").Value = 2090.6978
$ws.Range("J137").Value = 2622.3157
$ws.Range("L137").Value = 7866.9471
$ws.Range("N137").Value = -18066.9471

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 28000.8
$ws.Range("J19").Value = 30001
$ws.Range("L19").Value = 30001
$ws.Range("N19").Value = -30577

$ws.Range("H132").Value = 4823.387
$ws.Range("I132").Value = 6064.0625
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 18192.1875
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -15662.1875
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5860.8
$ws.Range("I40").Value = 6372.5713
$ws.Range("J40").Value = 4666.6665
$ws.Range("K40").Value = 6372.5713
$ws.Range("L40").Value = 4666.6665
$ws.Range("M40").Value = -6236.5713
$ws.Range("N40").Value = -4938.6665

$ws.Range("H45").Value = 7510.25
$ws.Range("I45").Value = 3041
$ws.Range("K45").Value = 3041
$ws.Range("M45").Value = -2634

$ws.Range("H122").Value = 2887.0205
$ws.Range("I122").Value = 2528
$ws.Range("J122").Value = 3784.5715
$ws.Range("K122").Value = 7584
$ws.Range("L122").Value = 11353.7145
$ws.Range("M122").Value = -5134
$ws.Range("N122").Value = -16253.7145

$ws.Range("H132").Value = 3197.8572
$ws.Range("I132").Value = 2364.2307
$ws.Range("K132").Value = 7092.6921
$ws.Range("M132").Value = -4562.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 1429067
$ws.Range("I3").Value = 5000001.5
$ws.Range("K3").Value = 5000001.5
$ws.Range("M3").Value = -4999887.5

$ws.Range("H107").Value = 737.35
$ws.Range("I107").Value = 252.92308
$ws.Range("J107").Value = 1637
$ws.Range("K107").Value = 758.76924
$ws.Range("L107").Value = 4911
$ws.Range("M107").Value = 1161.23076
$ws.Range("N107").Value = -8751

$ws.Range("H126").Value = 35418.234
$ws.Range("I126").Value = 43381.332
$ws.Range("J126").Value = 3565.8333
$ws.Range("K126").Value = 130143.996
$ws.Range("L126").Value = 10697.4999
$ws.Range("M126").Value = -127673.996
$ws.Range("N126").Value = -15637.4999

Write-Output "edits applied"
